$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.562.80"
$ws.Range("E2").Value = "  +2.37%  "

$ws.Range("D3").Value = "'1.665.80"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'237.23"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("E7").Value = "  -2.82%  "

$ws.Range("D8").Value = "'0.2577"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "'0.06133"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "'1.663.28"
$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").Value = "'0.06919"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("D13").Value = "'4.337"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").Value = "'75.16"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "'0.5716"
$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "'25.569.51"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").Value = "'0.000006675"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").Value = "'11.37"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("D21").Value = "'1.877.00"

$ws.Range("D22").Value = "'4.423"
$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("D23").Value = "'8.595"
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("D24").Value = "'5.218"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D25").Value = "'134.24"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "'14.88"
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").Value = "'1.376"
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("E28").Value = "  +5.38%  "

$ws.Range("D29").Value = "'104.07"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("D30").Value = "'3.941"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").Value = "'0.07654"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").Value = "'3.594"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("D33").Value = "'0.04328"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").Value = "'2.604"
$ws.Range("E34").Value = "  +1.22%  "

$ws.Range("D35").Value = "'0.6066"
$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("D36").Value = "'0.9389"
$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").Value = "'0.9263"
$ws.Range("E37").Value = "  +5.85%  "

$ws.Range("D38").Value = "'2.446"
$ws.Range("E38").Value = "  -5.31%  "

$ws.Range("D39").Value = "'106.79"
$ws.Range("E39").Value = "  +7.91%  "

$ws.Range("D40").Value = "'0.9996"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'1.830"
$ws.Range("E41").Value = "  +4.03%  "

$ws.Range("D42").Value = "'0.01449"
$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("D43").Value = "'5.066"
$ws.Range("E43").Value = "  +8.17%  "

$ws.Range("D44").Value = "'0.3705"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "'0.1110"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").Value = "'0.05264"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("D47").Value = "'6.085"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "'30.92"
$ws.Range("E48").Value = "  +7.45%  "

$ws.Range("D49").Value = "'7.574"
$ws.Range("E49").Value = "  +6.40%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'0.9990"
$ws.Range("E51").Value = "  +0.10%  "

